# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gains a new (blank) column before the
# existing "Late" column, shifting Late/heading/Outstanding one column to
# the right (N->O, O->P, P->Q). The active sheet/selection also changes
# from "Transactions"!D3 to "Repayment schedule"!M17.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Repayment schedule")
$ws4 = $wb.Worksheets.Item("Transactions")

# Insert a new blank column before column N; this shifts the existing
# N:P ("Late", "heading", "Outstanding") data right to O:Q.
$ws3.Columns("N").Insert()

# Give the freshly inserted column the same rendered width as its
# neighbours (~10.71 chars).
$ws3.Columns("N").ColumnWidth = 9.86

# The workbook was resaved with "Repayment schedule" as the active sheet
# (previously "Transactions" was active), with a new selected cell.
$ws3.Activate()
$ws3.Range("M17").Select() | Out-Null
